# Apply the "cryptos list" update (GitHub Actions commit, 2023-11-13).
# Updates price (column D) and 1h volume-change (column E) figures for most
# rows, and for rows 44/45 also swaps the Coin/Link data (HuobiToken now
# ranks above Aave, with updated price/volume figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $val) {
    # Assigning a numeric-looking string via .Value normally gets
    # auto-converted to a number by Excel. Force text storage (matching the
    # original inline-string cells) by temporarily switching the cell to the
    # "Text" number format, then restore the cell's original style so no
    # visible formatting/style change is introduced.
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$updates = @{
    'D2' = '36.754.42'
    'E2' = '  -1.24%  '
    'D3' = '2.090.63'
    'E3' = '  +1.70%  '
    'E4' = '  -0.02%  '
    'D5' = '245.24'
    'E5' = '  -1.44%  '
    'E6' = '  -1.90%  '
    'E7' = '  +0.05%  '
    'D8' = '54.14'
    'E8' = '  -6.31%  '
    'D9' = '58.81'
    'E9' = '  -2.18%  '
    'E10' = '  -4.36%  '
    'D11' = '0.0761'
    'E11' = '  -2.45%  '
    'E12' = '  +0.77%  '
    'D13' = '0.904'
    'E13' = '  +2.14%  '
    'D14' = '15.02'
    'E14' = '  -6.47%  '
    'D15' = '2.396.41'
    'E15' = '  +1.72%  '
    'D16' = '5.51'
    'E16' = '  -3.81%  '
    'D17' = '2.091.26'
    'E17' = '  +1.68%  '
    'D18' = '36.737.72'
    'E18' = '  -1.30%  '
    'D19' = '17.11'
    'E19' = '  -6.45%  '
    'D20' = '72.63'
    'E20' = '  -3.16%  '
    'D21' = '0.0₃0882'
    'E21' = '  -1.49%  '
    'D22' = '5.44'
    'E22' = '  +0.58%  '
    'D23' = '238.84'
    'E23' = '  +0.53%  '
    'E24' = '  +0.00%  '
    'E25' = '  -3.39%  '
    'D26' = '9.73'
    'E26' = '  +2.16%  '
    'D27' = '2.16'
    'E27' = '  -1.58%  '
    'D28' = '166.99'
    'E28' = '  -1.63%  '
    'E29' = '  +2.26%  '
    'E30' = '  -1.70%  '
    'D31' = '5.30'
    'E31' = '  +9.48%  '
    'D32' = '1.16'
    'E32' = '  +2.43%  '
    'D33' = '4.72'
    'E33' = '  +5.03%  '
    'D34' = '0.0609'
    'E34' = '  -1.87%  '
    'D35' = '2.42'
    'E35' = '  +7.31%  '
    'E36' = '  +0.19%  '
    'E37' = '  +3.92%  '
    'D38' = '0.0829'
    'E38' = '  -7.08%  '
    'E40' = '  +0.94%  '
    'D41' = '0.0220'
    'E41' = '  -1.39%  '
    'E42' = '  -7.86%  '
    'D43' = '0.0954'
    'E43' = '  -4.13%  '
    'B44' = 'HuobiToken'
    'C44' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D44' = '2.87'
    'E44' = '  -8.96%  '
    'B45' = 'Aave'
    'C45' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D45' = '96.19'
    'E45' = '  -0.19%  '
    'D46' = '16.06'
    'E46' = '  -7.47%  '
    'D47' = '1.384.85'
    'E47' = '  +8.80%  '
    'E48' = '  +8.47%  '
    'E49' = '  +0.08%  '
    'E50' = '  +1.37%  '
    'D51' = '2.283.04'
    'E51' = '  +1.75%  '
}

foreach ($cellRef in $updates.Keys) {
    Set-CellText $cellRef $updates[$cellRef]
}
